$d = $word.ActiveDocument

# Start from the last paragraph ("NASI POSTER") and append two more
# paragraphs after it, carrying forward the same run/paragraph
# formatting (bold, complex-script bold, 52pt / sz 104, centered,
# "Normal" style) -- the first new paragraph stays empty, the second
# gets the new text "Hope this will work".

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.InsertParagraphAfter()

$emptyPara = $d.Paragraphs.Last
$emptyPara.Range.InsertParagraphAfter()

$textPara = $d.Paragraphs.Last
$textPara.Range.Text = "Hope this will work"

Write-Output $d.Paragraphs.Count
